$wb = $excel.ActiveWorkbook

# --- Update the "Hoja1" sheet text in cell A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$rngA1 = $wsHoja1.Range("A1")
$text = $rngA1.Value()
$text = $text.Replace("✅ 1000 Bs = 12.77 = 51576.66 pesos", "✅ 1000 Bs = 12.79 = 51701.07 pesos")
$text = $text.Replace("✅ 51576.66 pesos = 12.68 = 979.03 Bs", "✅ 51701.07 pesos = 12.74 = 955.3 Bs")
$rngA1.Value = $text

# --- Update the "tasas" sheet numeric values ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 78.18000000000001
$wsTasas.Range("O10").Value = 4041.99
$wsTasas.Range("N12").Value = 4059
$wsTasas.Range("O12").Value = 75
